# "new Excel Model changes"
# Rename the "New Billing" label in A2 to "New_Billing", and move the
# active selection from D2 to B3 (matches saved sheetView selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "New_Billing"

$ws.Range("B3").Select()
